$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Move the old row 23 ("Pneumatics" / VRM channel) down to row 26,
#    then fill rows 23-25 with the three new Digital sensor pins.
# ---------------------------------------------------------------------

# -- old row23 content, re-homed at row 26 --
$ws.Range("B26").Value = "Pneumatics"
$ws.Range("C26").Value = 0
$ws.Range("D26").Formula = '=IF(EXACT(B26,"Analog"),IF(C26<4,"RoboRio","MXP"),IF(EXACT(B26,"Digital"),IF(C26<10,"RoboRio","MXP"),IF(EXACT(B26,"PWM"),IF(C26<10,"RoboRio","MXP"),"N/A")))'

# -- row 23: Tote Intake Sensor --
$ws.Range("A23").Value = "Tote Intake Sensor"
$ws.Range("B23").Value = "Digital"
$ws.Range("C23").Value = 7
$ws.Range("D23").Formula = '=IF(EXACT(B23,"Analog"),IF(C23<4,"RoboRio","MXP"),IF(EXACT(B23,"Digital"),IF(C23<10,"RoboRio","MXP"),IF(EXACT(B23,"PWM"),IF(C23<10,"RoboRio","MXP"),"N/A")))'

# -- row 24: Ejector Out --
$ws.Range("A24").Value = "Ejector Out"
$ws.Range("B24").Value = "Digital"
$ws.Range("C24").Value = 8
$ws.Range("D24").Formula = '=IF(EXACT(B24,"Analog"),IF(C24<4,"RoboRio","MXP"),IF(EXACT(B24,"Digital"),IF(C24<10,"RoboRio","MXP"),IF(EXACT(B24,"PWM"),IF(C24<10,"RoboRio","MXP"),"N/A")))'

# -- row 25: Ejector In --
$ws.Range("A25").Value = "Ejector In"
$ws.Range("B25").Value = "Digital"
$ws.Range("C25").Value = 9
$ws.Range("D25").Formula = '=IF(EXACT(B25,"Analog"),IF(C25<4,"RoboRio","MXP"),IF(EXACT(B25,"Digital"),IF(C25<10,"RoboRio","MXP"),IF(EXACT(B25,"PWM"),IF(C25<10,"RoboRio","MXP"),"N/A")))'

# carry the "styled group" formatting (s="1") down across A23:C26, matching
# the banded look the rest of the lower table rows already use
$ws.Range("A18:C26").Style = "Normal"

# give the new calculated-column cells a distinct (border-less) number
# format so they pick up their own style slot, same as the source file
$ws.Range("D23:D25").NumberFormat = "General"

# C22 (Lift Min/Max Switch pin) bumped from 2 to 6
$ws.Range("C22").Value = 6

# ---------------------------------------------------------------------
# 2. Grow the table to cover the new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D26"))

# ---------------------------------------------------------------------
# 3. Column width tweaks.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666
$ws.Columns.Item(3).ColumnWidth = 4.25

# ---------------------------------------------------------------------
# 4. View tweaks - zoom + active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("C23").Select()
